$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.868.72"
$ws.Range("E2").Value = "  +4.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.274.88"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.93"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.30"
$ws.Range("E6").Value = "  +6.37%  "
$ws.Range("E7").Value = "  +3.94%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +3.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.70"
$ws.Range("E10").Value = "  +6.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.45"
$ws.Range("E11").Value = "  +4.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.115"
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("E14").Value = "  +3.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.624.66"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("E16").Value = "  +3.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.290.43"
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("E18").Value = "  +3.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.792.46"
$ws.Range("E19").Value = "  +4.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.25"
$ws.Range("E20").Value = "  +8.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0910"
$ws.Range("E21").Value = "  +2.51%  "
$ws.Range("E22").Value = "  +3.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.32"
$ws.Range("E23").Value = "  +2.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "243.79"
$ws.Range("E24").Value = "  +3.01%  "
$ws.Range("E25").Value = "  +3.80%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  +5.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.29"
$ws.Range("E28").Value = "  +4.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.62"
$ws.Range("E29").Value = "  +3.33%  "
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("E31").Value = "  +7.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.49"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  +4.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0752"
$ws.Range("E35").Value = "  +4.94%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  +3.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.72"
$ws.Range("E38").Value = "  +8.25%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.116"
$ws.Range("E39").Value = "  +2.92%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.105"
$ws.Range("E40").Value = "  +5.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.83"
$ws.Range("E41").Value = "  +4.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.93"
$ws.Range("E42").Value = "  +5.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.072.17"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.60"
$ws.Range("E44").Value = "  +5.58%  "
$ws.Range("E45").Value = "  +3.16%  "
$ws.Range("E46").Value = "  +3.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.91"
$ws.Range("E47").Value = "  +7.07%  "
$ws.Range("E48").Value = "  +4.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.53"
$ws.Range("E49").Value = "  +3.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.06"
$ws.Range("E50").Value = "  +7.52%  "
$ws.Range("E51").Value = "  +3.39%  "
